$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.368.57'
$ws.Range('E2').Value = '  +0.12%  '

$ws.Range('D3').Value = '1.872.28'
$ws.Range('E3').Value = '  -0.62%  '

$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '235.51'
$ws.Range('E5').Value = '  -1.10%  '

$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.02%  '

$ws.Range('D7').Value = '0.4671'
$ws.Range('E7').Value = '  -0.23%  '

$ws.Range('D8').Value = '0.2843'
$ws.Range('E8').Value = '  +0.53%  '

$ws.Range('D9').Value = '0.06566'
$ws.Range('E9').Value = '  -0.07%  '

$ws.Range('D10').Value = '21.36'
$ws.Range('E10').Value = '  +7.82%  '

$ws.Range('D11').Value = '0.07921'
$ws.Range('E11').Value = '  +1.97%  '

$ws.Range('D12').Value = '97.51'
$ws.Range('E12').Value = '  -0.65%  '

$ws.Range('D13').Value = '1.863.07'
$ws.Range('E13').Value = '  -1.11%  '

$ws.Range('D14').Value = '5.135'
$ws.Range('E14').Value = '  +0.02%  '

$ws.Range('D15').Value = '0.6765'
$ws.Range('E15').Value = '  +1.35%  '

$ws.Range('D16').Value = '281.25'
$ws.Range('E16').Value = '  -1.22%  '

$ws.Range('D17').Value = '30.359.60'
$ws.Range('E17').Value = '  +0.07%  '

$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.10%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '5.522'
$ws.Range('E19').Value = '  +3.10%  '

$ws.Range('D20').Value = '12.68'
$ws.Range('E20').Value = '  +0.47%  '

$ws.Range('D21').Value = '2.115.15'
$ws.Range('E21').Value = '  -0.68%  '

$ws.Range('D22').Value = '0.000007301'
$ws.Range('E22').Value = '  -0.21%  '

$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('D24').Value = '6.206'
$ws.Range('E24').Value = '  +0.54%  '

$ws.Range('D25').Value = '9.277'
$ws.Range('E25').Value = '  +0.05%  '

$ws.Range('D26').Value = '165.04'
$ws.Range('E26').Value = '  -1.53%  '

$ws.Range('D27').Value = '19.16'
$ws.Range('E27').Value = '  +0.40%  '

$ws.Range('D28').Value = '1.944'
$ws.Range('E28').Value = '  -2.66%  '

$ws.Range('D29').Value = '1.375'
$ws.Range('E29').Value = '  -0.07%  '

$ws.Range('D30').Value = '0.09728'
$ws.Range('E30').Value = '  -1.13%  '

$ws.Range('D31').Value = '4.439'
$ws.Range('E31').Value = '  -0.65%  '

$ws.Range('E32').Value = '  -1.20%  '

$ws.Range('D33').Value = '4.117'
$ws.Range('E33').Value = '  -1.60%  '

$ws.Range('D34').Value = '0.04689'
$ws.Range('E34').Value = '  -0.53%  '

$ws.Range('D35').Value = '1.118'
$ws.Range('E35').Value = '  +2.00%  '

# Leading "'" forces text so the trailing zero in "0.7060" is preserved
# (otherwise Excel would coerce it to the number 0.706).
$ws.Range('D36').Value = '''0.7060'
$ws.Range('E36').Value = '  -0.55%  '

$ws.Range('E37').Value = '  -0.22%  '

$ws.Range('E38').Value = '  -0.46%  '

$ws.Range('D39').Value = '6.327'
$ws.Range('E39').Value = '  -6.08%  '

$ws.Range('D40').Value = '2.542'
$ws.Range('E40').Value = '  +0.84%  '

$ws.Range('D41').Value = '73.29'
$ws.Range('E41').Value = '  +1.32%  '

$ws.Range('D42').Value = '1.944'
$ws.Range('E42').Value = '  -1.59%  '

$ws.Range('D43').Value = '0.8492'
$ws.Range('E43').Value = '  -2.62%  '

$ws.Range('D44').Value = '0.4191'
$ws.Range('E44').Value = '  -0.33%  '

$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.09%  '

$ws.Range('D46').Value = '103.87'
$ws.Range('E46').Value = '  -0.39%  '

$ws.Range('D47').Value = '7.228'
$ws.Range('E47').Value = '  -0.20%  '

$ws.Range('D48').Value = '9.149'
$ws.Range('E48').Value = '  -2.17%  '

$ws.Range('D49').Value = '932.55'
$ws.Range('E49').Value = '  -6.00%  '

$ws.Range('D50').Value = '34.16'
$ws.Range('E50').Value = '  +0.00%  '

$ws.Range('E51').Value = '  -2.60%  '
